$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.754434108734131
$ws.Range("B1").Value = 1.875431060791016
$ws.Range("C1").Value = 5.045635223388672
$ws.Range("D1").Value = 1.745492219924927
$ws.Range("E1").Value = 0.4892308712005615
